$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh crypto "Price" (D) and "Volume(1h)" (E) columns per the
# Oct 14 2023 GitHub Actions data pull. Force text format so values
# such as "217.20" or "26.911.77" are not coerced into numbers,
# matching the original inlineStr text cells.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.911.77"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.29%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.551.85"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.29%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.62%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "206.54"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.38%  "
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.08%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.62%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "22.02"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +1.58%  "
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -0.45%  "
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +0.73%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0856"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -0.72%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.773.46"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -0.29%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.540.57"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -1.11%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.76"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +0.90%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.519"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +0.65%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "26.906.07"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -0.34%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "61.62"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.61%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0₃0710"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +2.99%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "217.20"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.57%  "
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.38%  "
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.61%  "
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +1.01%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.18"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.12%  "
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -1.45%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "153.40"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +0.66%  "
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.20%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "14.99"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +0.47%  "
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +0.68%  "
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.58%  "
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +1.42%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.09"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -0.82%  "
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -0.38%  "
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +3.56%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.411.09"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +0.35%  "
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +2.05%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.968"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +1.09%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.30"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +0.25%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0165"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -0.07%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.528"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +1.02%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.807"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -0.39%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.66"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +2.57%  "
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +1.48%  "
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +0.70%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "64.55"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +0.90%  "
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -1.30%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.687.06"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -0.33%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "87.14"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +1.00%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0521"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +1.75%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0₆01000"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +1.94%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0959"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +0.22%  "
